$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 84.77257745550429
$ws.Range("D2").Value = 83.52649006622516
$ws.Range("E2").Value = 79.32389937106919
$ws.Range("F2").Value = 81.37096774193549

$ws.Range("C3").Value = 85.33289386947924
$ws.Range("D3").Value = 84.48707256046706
$ws.Range("E3").Value = 79.63836477987421
$ws.Range("F3").Value = 81.99109672197491

$ws.Range("C4").Value = 84.64073829927489
$ws.Range("D4").Value = 84.13910093299405
$ws.Range("E4").Value = 78.04878048780488
$ws.Range("F4").Value = 80.9795918367347

$ws.Range("C5").Value = 85.43177323665128
$ws.Range("D5").Value = 84.45552784704904
$ws.Range("E5").Value = 79.93705743509048
$ws.Range("F5").Value = 82.13419563459983

$ws.Range("C6").Value = 85.03132212331025
$ws.Range("D6").Value = 84.35660218671153
$ws.Range("E6").Value = 78.91424075531079
$ws.Range("F6").Value = 81.54471544715449
